# Swap the content of rows 8 and 9 for the columns that differ between the
# two fungi records (A, B, E, F, G, H, I, P, Q, R). The other columns are
# identical between the two rows so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "I", "P", "Q", "R")

foreach ($col in $cols) {
    $cell8 = $ws.Range("$col`8")
    $cell9 = $ws.Range("$col`9")

    $v8 = $cell8.Value2
    $v9 = $cell9.Value2

    $cell8.Value2 = $v9
    $cell9.Value2 = $v8
}
